$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 8 (shifts rows 8..42 down to 9..43)
$ws.Rows.Item(8).Insert()

# Copy the formatting of the "Index3" header row (row 7) into the new row 8
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)

# Populate the new "Index4" row
$ws.Range("A8").Value() = "Index4"
$ws.Range("C8").Value() = "AcctCode,CustNo,FacmNo,ClsFlag,RvNo"

# Update the view: scroll back to the top-left and move the selection to C7
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
